# Auto-generated edit script applying the scheduled-runner market-data refresh
# to the Ragnarok_Profits workbook (H/I/J/K/L/M/N columns per crafting-job sheet).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 1476.4286
$ws.Cells.Item(4, 9).Value = 1380.5714
$ws.Cells.Item(4, 11).Value = 1380.5714
$ws.Cells.Item(4, 13).Value = -1266.5714
$ws.Cells.Item(5, 8).Value = 68.63636
$ws.Cells.Item(5, 9).Value = 77.57143000000001
$ws.Cells.Item(5, 11).Value = 77.57143000000001
$ws.Cells.Item(5, 13).Value = 37.42856999999999
$ws.Cells.Item(33, 8).Value = 284.8
$ws.Cells.Item(33, 9).Value = 258.5
$ws.Cells.Item(33, 11).Value = 258.5
$ws.Cells.Item(33, 13).Value = -29.5
$ws.Cells.Item(64, 8).Value = 4031.5454
$ws.Cells.Item(64, 10).Value = 4059.75
$ws.Cells.Item(64, 12).Value = 4059.75
$ws.Cells.Item(64, 14).Value = -4555.75
$ws.Cells.Item(67, 8).Value = 4031.5454
$ws.Cells.Item(67, 10).Value = 4059.75
$ws.Cells.Item(67, 12).Value = 4059.75
$ws.Cells.Item(67, 14).Value = -5775.75
$ws.Cells.Item(88, 8).Value = 2702.389
$ws.Cells.Item(88, 10).Value = 2578.25
$ws.Cells.Item(88, 12).Value = 2578.25
$ws.Cells.Item(88, 14).Value = -3390.25
$ws.Cells.Item(91, 8).Value = 2702.389
$ws.Cells.Item(91, 10).Value = 2578.25
$ws.Cells.Item(91, 12).Value = 2578.25
$ws.Cells.Item(91, 14).Value = -5386.25
$ws.Cells.Item(106, 8).Value = 7524.7856
$ws.Cells.Item(106, 9).Value = 7263.364
$ws.Cells.Item(106, 11).Value = 7263.364
$ws.Cells.Item(106, 13).Value = -6632.364
$ws.Cells.Item(113, 8).Value = 7166.6665
$ws.Cells.Item(113, 9).Value = 7250
$ws.Cells.Item(113, 10).Value = 7000
$ws.Cells.Item(113, 11).Value = 7250
$ws.Cells.Item(113, 12).Value = 7000
$ws.Cells.Item(113, 13).Value = -3996
$ws.Cells.Item(113, 14).Value = -13508
$ws.Cells.Item(135, 8).Value = 3726.1667
$ws.Cells.Item(135, 9).Value = 832.1429000000001
$ws.Cells.Item(135, 10).Value = 7777.8
$ws.Cells.Item(135, 11).Value = 7489.2861
$ws.Cells.Item(135, 12).Value = 70000.2
$ws.Cells.Item(135, 13).Value = -4954.2861
$ws.Cells.Item(135, 14).Value = -75070.2
$ws.Cells.Item(137, 8).Value = 2115.3914
$ws.Cells.Item(137, 9).Value = 2042.7368
$ws.Cells.Item(137, 11).Value = 6128.2104
$ws.Cells.Item(137, 13).Value = -3578.2104

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3094.2917
$ws.Cells.Item(2, 9).Value = 3046.2
$ws.Cells.Item(2, 10).Value = 3174.4443
$ws.Cells.Item(2, 11).Value = 3046.2
$ws.Cells.Item(2, 12).Value = 3174.4443
$ws.Cells.Item(2, 13).Value = -2933.2
$ws.Cells.Item(2, 14).Value = -3400.4443
$ws.Cells.Item(43, 8).Value = 37866.11
$ws.Cells.Item(43, 10).Value = 37486.43
$ws.Cells.Item(43, 12).Value = 37486.43
$ws.Cells.Item(43, 14).Value = -38112.43
$ws.Cells.Item(61, 8).Value = 16430605
$ws.Cells.Item(61, 9).Value = 23335164
$ws.Cells.Item(61, 11).Value = 23335164
$ws.Cells.Item(61, 13).Value = -23334952
$ws.Cells.Item(74, 8).Value = 2690.2778
$ws.Cells.Item(74, 9).Value = 1537.6428
$ws.Cells.Item(74, 10).Value = 6724.5
$ws.Cells.Item(74, 11).Value = 1537.6428
$ws.Cells.Item(74, 12).Value = 6724.5
$ws.Cells.Item(74, 13).Value = -663.6428000000001
$ws.Cells.Item(74, 14).Value = -8472.5
$ws.Cells.Item(77, 8).Value = 2690.2778
$ws.Cells.Item(77, 9).Value = 1537.6428
$ws.Cells.Item(77, 10).Value = 6724.5
$ws.Cells.Item(77, 11).Value = 7688.214
$ws.Cells.Item(77, 12).Value = 33622.5
$ws.Cells.Item(77, 13).Value = -3320.214
$ws.Cells.Item(77, 14).Value = -42358.5
$ws.Cells.Item(97, 8).Value = 2201
$ws.Cells.Item(97, 9).Value = 2246.7778
$ws.Cells.Item(97, 11).Value = 2246.7778
$ws.Cells.Item(97, 13).Value = -1750.7778
$ws.Cells.Item(110, 8).Value = 2949.5
$ws.Cells.Item(110, 9).Value = 899.5
$ws.Cells.Item(110, 11).Value = 899.5
$ws.Cells.Item(110, 13).Value = 1145.5
$ws.Cells.Item(116, 8).Value = 3094.2917
$ws.Cells.Item(116, 9).Value = 3046.2
$ws.Cells.Item(116, 10).Value = 3174.4443
$ws.Cells.Item(116, 11).Value = 3046.2
$ws.Cells.Item(116, 12).Value = 3174.4443
$ws.Cells.Item(116, 13).Value = -752.1999999999998
$ws.Cells.Item(116, 14).Value = -7762.4443
$ws.Cells.Item(132, 8).Value = 2783416.2
$ws.Cells.Item(132, 9).Value = 5334.793
$ws.Cells.Item(132, 11).Value = 16004.379
$ws.Cells.Item(132, 13).Value = -13474.379
$ws.Cells.Item(136, 8).Value = 16430605
$ws.Cells.Item(136, 9).Value = 23335164
$ws.Cells.Item(136, 11).Value = 70005492
$ws.Cells.Item(136, 13).Value = -70002942

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3094.2917
$ws.Cells.Item(3, 9).Value = 3046.2
$ws.Cells.Item(3, 10).Value = 3174.4443
$ws.Cells.Item(3, 11).Value = 3046.2
$ws.Cells.Item(3, 12).Value = 3174.4443
$ws.Cells.Item(3, 13).Value = -2932.2
$ws.Cells.Item(3, 14).Value = -3402.4443
$ws.Cells.Item(134, 8).Value = 4350128.5
$ws.Cells.Item(134, 9).Value = 2248.5
$ws.Cells.Item(134, 11).Value = 6745.5
$ws.Cells.Item(134, 13).Value = -4210.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3000
$ws.Cells.Item(62, 9).Value = 1505
$ws.Cells.Item(62, 10).Value = 3373.75
$ws.Cells.Item(62, 11).Value = 1505
$ws.Cells.Item(62, 12).Value = 3373.75
$ws.Cells.Item(62, 13).Value = -881
$ws.Cells.Item(62, 14).Value = -4621.75
$ws.Cells.Item(65, 8).Value = 3000
$ws.Cells.Item(65, 9).Value = 1505
$ws.Cells.Item(65, 10).Value = 3373.75
$ws.Cells.Item(65, 11).Value = 7525
$ws.Cells.Item(65, 12).Value = 16868.75
$ws.Cells.Item(65, 13).Value = -4405
$ws.Cells.Item(65, 14).Value = -23108.75
$ws.Cells.Item(82, 8).Value = 99998.164
$ws.Cells.Item(82, 10).Value = 99998.164
$ws.Cells.Item(82, 12).Value = 99998.164
$ws.Cells.Item(82, 14).Value = -100720.164
$ws.Cells.Item(85, 8).Value = 99998.164
$ws.Cells.Item(85, 10).Value = 99998.164
$ws.Cells.Item(85, 12).Value = 99998.164
$ws.Cells.Item(85, 14).Value = -102494.164
$ws.Cells.Item(107, 8).Value = 960.26666
$ws.Cells.Item(107, 9).Value = 644.6923
$ws.Cells.Item(107, 10).Value = 3011.5
$ws.Cells.Item(107, 11).Value = 644.6923
$ws.Cells.Item(107, 12).Value = 3011.5
$ws.Cells.Item(107, 13).Value = 1275.3077
$ws.Cells.Item(107, 14).Value = -6851.5
$ws.Cells.Item(116, 8).Value = 99999
$ws.Cells.Item(116, 10).Value = 99999
$ws.Cells.Item(116, 12).Value = 99999
$ws.Cells.Item(116, 14).Value = -109177
$ws.Cells.Item(132, 8).Value = 2854.4
$ws.Cells.Item(132, 10).Value = 2975.2856
$ws.Cells.Item(132, 12).Value = 8925.856800000001
$ws.Cells.Item(132, 14).Value = -13985.8568
$ws.Cells.Item(134, 8).Value = 3512.25
$ws.Cells.Item(134, 9).Value = 3585.4285
$ws.Cells.Item(134, 11).Value = 10756.2855
$ws.Cells.Item(134, 13).Value = -8221.2855

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 7942.4
$ws.Cells.Item(9, 9).Value = 1793
$ws.Cells.Item(9, 11).Value = 5379
$ws.Cells.Item(9, 13).Value = -5155
$ws.Cells.Item(15, 8).Value = 6809.6
$ws.Cells.Item(15, 10).Value = 11316
$ws.Cells.Item(15, 12).Value = 33948
$ws.Cells.Item(15, 14).Value = -34228
$ws.Cells.Item(23, 8).Value = 1215.2354
$ws.Cells.Item(23, 9).Value = 67.71429000000001
$ws.Cells.Item(23, 10).Value = 2018.5
$ws.Cells.Item(23, 11).Value = 203.14287
$ws.Cells.Item(23, 12).Value = 6055.5
$ws.Cells.Item(23, 13).Value = 31.85712999999998
$ws.Cells.Item(23, 14).Value = -6525.5
$ws.Cells.Item(24, 8).Value = 9888.75
$ws.Cells.Item(24, 10).Value = 16777.5
$ws.Cells.Item(24, 12).Value = 50332.5
$ws.Cells.Item(24, 14).Value = -50792.5
$ws.Cells.Item(58, 8).Value = 10610.75
$ws.Cells.Item(58, 10).Value = 19666.5
$ws.Cells.Item(58, 12).Value = 58999.5
$ws.Cells.Item(58, 14).Value = -59255.5
$ws.Cells.Item(86, 8).Value = 590.3333
$ws.Cells.Item(86, 10).Value = 852.1667
$ws.Cells.Item(86, 12).Value = 2556.5001
$ws.Cells.Item(86, 14).Value = -4928.5001
$ws.Cells.Item(89, 8).Value = 590.3333
$ws.Cells.Item(89, 10).Value = 852.1667
$ws.Cells.Item(89, 12).Value = 7669.5003
$ws.Cells.Item(89, 14).Value = -19525.5003
$ws.Cells.Item(107, 8).Value = 11376582
$ws.Cells.Item(107, 9).Value = 290
$ws.Cells.Item(107, 10).Value = 15168679
$ws.Cells.Item(107, 11).Value = 870
$ws.Cells.Item(107, 12).Value = 45506037
$ws.Cells.Item(107, 13).Value = 1050
$ws.Cells.Item(107, 14).Value = -45509877

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2060667.4
$ws.Cells.Item(113, 9).Value = 3399
$ws.Cells.Item(113, 11).Value = 3399
$ws.Cells.Item(113, 13).Value = -1229
$ws.Cells.Item(132, 8).Value = 6671547
$ws.Cells.Item(132, 9).Value = 4664
$ws.Cells.Item(132, 11).Value = 13992
$ws.Cells.Item(132, 13).Value = -11462

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 20699.8
$ws.Cells.Item(22, 9).Value = 25499.75
$ws.Cells.Item(22, 11).Value = 25499.75
$ws.Cells.Item(22, 13).Value = -25204.75
$ws.Cells.Item(27, 8).Value = 20699.8
$ws.Cells.Item(27, 9).Value = 25499.75
$ws.Cells.Item(27, 11).Value = 25499.75
$ws.Cells.Item(27, 13).Value = -25392.75
$ws.Cells.Item(61, 8).Value = 8449
$ws.Cells.Item(61, 9).Value = 1899.75
$ws.Cells.Item(61, 11).Value = 1899.75
$ws.Cells.Item(61, 13).Value = -1697.75
$ws.Cells.Item(113, 8).Value = 8449
$ws.Cells.Item(113, 9).Value = 1899.75
$ws.Cells.Item(113, 11).Value = 1899.75
$ws.Cells.Item(113, 13).Value = 270.25
$ws.Cells.Item(122, 8).Value = 3736.8125
$ws.Cells.Item(122, 9).Value = 3335.72
$ws.Cells.Item(122, 10).Value = 5169.2856
$ws.Cells.Item(122, 11).Value = 10007.16
$ws.Cells.Item(122, 12).Value = 15507.8568
$ws.Cells.Item(122, 13).Value = -7557.16
$ws.Cells.Item(122, 14).Value = -20407.8568

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 21495.5
$ws.Cells.Item(41, 10).Value = 20395.6
$ws.Cells.Item(41, 12).Value = 20395.6
$ws.Cells.Item(41, 14).Value = -21175.6
$ws.Cells.Item(113, 8).Value = 689.9048
$ws.Cells.Item(113, 9).Value = 577.25
$ws.Cells.Item(113, 11).Value = 1731.75
$ws.Cells.Item(113, 13).Value = 438.25
$ws.Cells.Item(117, 8).Value = 99999
$ws.Cells.Item(117, 10).Value = 99999
$ws.Cells.Item(117, 12).Value = 99999
$ws.Cells.Item(117, 14).Value = -109177
$ws.Cells.Item(132, 8).Value = 424991.88
$ws.Cells.Item(132, 9).Value = 7937.3687
$ws.Cells.Item(132, 11).Value = 23812.1061
$ws.Cells.Item(132, 13).Value = -21282.1061
